$wb = $excel.ActiveWorkbook

# ALC row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("K8").Value = 300
$ws.Range("M8").Value = -161

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5638.815
$ws.Range("I51").Value = 4071.4285
$ws.Range("J51").Value = 7326.769
$ws.Range("K51").Value = 4071.4285
$ws.Range("L51").Value = 7326.769
$ws.Range("M51").Value = -3587.4285
$ws.Range("N51").Value = -8294.769

# ALC row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3456.8333
$ws.Range("I98").Value = 2949.25
$ws.Range("K98").Value = 2949.25
$ws.Range("M98").Value = -1451.25

# ALC row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 3456.8333
$ws.Range("I122").Value = 2949.25
$ws.Range("K122").Value = 8847.75
$ws.Range("M122").Value = -6397.75

# ALC row 131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 3761.5715
$ws.Range("I131").Value = 785.3570999999999
$ws.Range("J131").Value = 9714
$ws.Range("K131").Value = 2356.0713
$ws.Range("L131").Value = 29142
$ws.Range("M131").Value = 2683.9287
$ws.Range("N131").Value = -39222

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 115805.07
$ws.Range("I132").Value = 292977.12
$ws.Range("K132").Value = 878931.36
$ws.Range("M132").Value = -876401.36

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 4980.125
$ws.Range("I137").Value = 1787.25
$ws.Range("K137").Value = 5361.75
$ws.Range("M137").Value = -2811.75

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 7131.7856
$ws.Range("J138").Value = 7996.857
$ws.Range("L138").Value = 23990.571
$ws.Range("N138").Value = -34270.571

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 819428.75
$ws.Range("I2").Value = 1345207.5
$ws.Range("K2").Value = 1345207.5
$ws.Range("M2").Value = -1345094.5

# ARM row 4
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 289.2857
$ws.Range("I4").Value = 287.5
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 287.5
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = -171.5
$ws.Range("N4").Value = -532

# ARM row 5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 916
$ws.Range("I5").Value = 916
$ws.Range("K5").Value = 916
$ws.Range("M5").Value = -804

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3795.3035
$ws.Range("I32").Value = 2019.0667
$ws.Range("K32").Value = 2019.0667
$ws.Range("M32").Value = -1732.0667

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2004.8334
$ws.Range("J45").Value = 148
$ws.Range("L45").Value = 148
$ws.Range("N45").Value = -902

# ARM row 50
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 3499.75
$ws.Range("I50").Value = 5000
$ws.Range("J50").Value = 2999.6667
$ws.Range("K50").Value = 5000
$ws.Range("L50").Value = 2999.6667
$ws.Range("M50").Value = -4286
$ws.Range("N50").Value = -4427.6667

# ARM row 51
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 74749
$ws.Range("J51").Value = 74749
$ws.Range("L51").Value = 74749
$ws.Range("N51").Value = -76261

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1460300.1
$ws.Range("I110").Value = 1856291
$ws.Range("K110").Value = 1856291
$ws.Range("M110").Value = -1854246

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 819428.75
$ws.Range("I116").Value = 1345207.5
$ws.Range("K116").Value = 1345207.5
$ws.Range("M116").Value = -1342913.5

# ARM row 133
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 71606.60000000001
$ws.Range("J133").Value = 71606.60000000001
$ws.Range("L133").Value = 71606.60000000001
$ws.Range("N133").Value = -76666.60000000001

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 819428.75
$ws.Range("I3").Value = 1345207.5
$ws.Range("K3").Value = 1345207.5
$ws.Range("M3").Value = -1345093.5

# BSM row 4
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 916
$ws.Range("I4").Value = 916
$ws.Range("K4").Value = 916
$ws.Range("M4").Value = -801

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2977273
$ws.Range("I99").Value = 4167382.2
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 4167382.2
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -4165884.2
$ws.Range("N99").Value = -4996

# BSM row 140
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 144821.73
$ws.Range("J140").Value = 144821.73
$ws.Range("L140").Value = 144821.73
$ws.Range("N140").Value = -155181.73

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1733.25
$ws.Range("I16").Value = 1709.4546
$ws.Range("J16").Value = 1995
$ws.Range("K16").Value = 1709.4546
$ws.Range("L16").Value = 1995
$ws.Range("M16").Value = -1422.4546
$ws.Range("N16").Value = -2569

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4016.7036
$ws.Range("I31").Value = 987.375
$ws.Range("K31").Value = 987.375
$ws.Range("M31").Value = -692.375

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4016.7036
$ws.Range("I34").Value = 987.375
$ws.Range("K34").Value = 987.375
$ws.Range("M34").Value = -785.375

# CRP row 107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1013324.6
$ws.Range("I107").Value = 1516320.4
$ws.Range("J107").Value = 7333.1665
$ws.Range("K107").Value = 1516320.4
$ws.Range("L107").Value = 7333.1665
$ws.Range("M107").Value = -1514400.4
$ws.Range("N107").Value = -11173.1665

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 1733.25
$ws.Range("I113").Value = 1709.4546
$ws.Range("J113").Value = 1995
$ws.Range("K113").Value = 1709.4546
$ws.Range("L113").Value = 1995
$ws.Range("M113").Value = 460.5454
$ws.Range("N113").Value = -6335

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 9871.333000000001
$ws.Range("I132").Value = 4037
$ws.Range("J132").Value = 30291.5
$ws.Range("K132").Value = 12111
$ws.Range("L132").Value = 90874.5
$ws.Range("M132").Value = -9581
$ws.Range("N132").Value = -95934.5

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1557496.2
$ws.Range("J5").Value = 2593825.8
$ws.Range("L5").Value = 7781477.399999999
$ws.Range("N5").Value = -7781701.399999999

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 244.31818
$ws.Range("J12").Value = 220.875
$ws.Range("L12").Value = 662.625
$ws.Range("N12").Value = -1008.625

# CUL row 26
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1025.5
$ws.Range("I26").Value = 367.33334
$ws.Range("J26").Value = 3000
$ws.Range("K26").Value = 1102.00002
$ws.Range("L26").Value = 9000
$ws.Range("M26").Value = -814.0000199999999
$ws.Range("N26").Value = -9576

# CUL row 56
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6688.95
$ws.Range("I56").Value = 6688.95
$ws.Range("K56").Value = 6688.95
$ws.Range("M56").Value = -6158.95

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7694180.5
$ws.Range("J131").Value = 6412410.5
$ws.Range("L131").Value = 19237231.5
$ws.Range("N131").Value = -19247311.5

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1557496.2
$ws.Range("J135").Value = 2593825.8
$ws.Range("L135").Value = 23344432.2
$ws.Range("N135").Value = -23349502.2

# CUL row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 671667.6
$ws.Range("I136").Value = 771154.9399999999
$ws.Range("K136").Value = 2313464.82
$ws.Range("M136").Value = -2308364.82

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 15877016
$ws.Range("I70").Value = 47619050
$ws.Range("K70").Value = 47619050
$ws.Range("M70").Value = -47618780

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 15877016
$ws.Range("I73").Value = 47619050
$ws.Range("K73").Value = 47619050
$ws.Range("M73").Value = -47618114

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7575.643
$ws.Range("I102").Value = 7328.6816
$ws.Range("K102").Value = 7328.6816
$ws.Range("M102").Value = -5706.6816

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 9524662
$ws.Range("I107").Value = 23810024
$ws.Range("J107").Value = 1087.6666
$ws.Range("K107").Value = 23810024
$ws.Range("L107").Value = 1087.6666
$ws.Range("M107").Value = -23808104
$ws.Range("N107").Value = -4927.6666

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2513
$ws.Range("I132").Value = 2693.5625
$ws.Range("K132").Value = 8080.6875
$ws.Range("M132").Value = -5550.6875

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 884.3913
$ws.Range("I22").Value = 561
$ws.Range("K22").Value = 561
$ws.Range("M22").Value = -266

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 884.3913
$ws.Range("I27").Value = 561
$ws.Range("K27").Value = 561
$ws.Range("M27").Value = -454

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5946.926
$ws.Range("I46").Value = 3098.111
$ws.Range("K46").Value = 3098.111
$ws.Range("M46").Value = -2910.111

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5141.1665
$ws.Range("I61").Value = 6237.25
$ws.Range("K61").Value = 6237.25
$ws.Range("M61").Value = -6035.25

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1134.0667
$ws.Range("I100").Value = 1000.9231
$ws.Range("K100").Value = 1000.9231
$ws.Range("M100").Value = -459.9231

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 5141.1665
$ws.Range("I113").Value = 6237.25
$ws.Range("K113").Value = 6237.25
$ws.Range("M113").Value = -4067.25

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5401.3447
$ws.Range("I132").Value = 3332.2307
$ws.Range("J132").Value = 7082.5
$ws.Range("K132").Value = 9996.6921
$ws.Range("L132").Value = 21247.5
$ws.Range("M132").Value = -7466.6921
$ws.Range("N132").Value = -26307.5

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5244.4
$ws.Range("I136").Value = 4055.5
$ws.Range("K136").Value = 12166.5
$ws.Range("M136").Value = -9616.5
